# Adds "accuracy_balanced_mean" (new rows 16-22) and "accuracy_balanced_std"
# (new rows 37-43) metric blocks to the metrics table, per reviewer feedback.
# Existing f1_macro_std / f1_micro_std blocks shift down by 7 rows but keep
# their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert 7 blank rows above the current "f1_macro_std" block (row 16)
#    so the new "accuracy_balanced_mean" block can be written there.
# ---------------------------------------------------------------------
$ws.Range("A16:A22").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2) Populate the new "accuracy_balanced_mean" rows (16-22).
# ---------------------------------------------------------------------
$accBalMeanRows = @(
    @{ Row=16; N="0";     C=0;     D=0;     E=0;     F=0;     G=0;     H=0.437 },
    @{ Row=17; N="100";   C=0.185; D=0.167; E=0.287; F=0.254; G=0.306; H=0.527 },
    @{ Row=18; N="500";   C=0.404; D=0.379; E=0.422; F=0.442; G=0.598; H=0.617 },
    @{ Row=19; N="1000";  C=0.454; D=0.44;  E=0.495; F=0.483; G=0.652; H=0.637 },
    @{ Row=20; N="2500";  C=0.501; D=0.51;  E=0.547; F=0.553; G=0.694; H=0.673 },
    @{ Row=21; N="5000";  C=0.542; D=0.548; E=0.586; F=0.601; G=0.71;  H=0.7   },
    @{ Row=22; N="10000"; C=0.5679999999999999; D=0.587; E=0.609; F=0.629; G=0.741; H=0.719 }
)

foreach ($r in $accBalMeanRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "accuracy_balanced_mean"
    $ws.Range("B$row").Value = "'" + $r.N
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
}

# ---------------------------------------------------------------------
# 3) Append 7 new rows at the bottom (37-43) for "accuracy_balanced_std".
#    The sheet currently ends at row 36 after the insert above, so these
#    rows are brand new, not an insert-and-shift.
# ---------------------------------------------------------------------
$accBalStdRows = @(
    @{ Row=37; N="0";     C=0;     D=0;     E=0;     F=0;     G=0;     H=0     },
    @{ Row=38; N="100";   C=0.011; D=0.004; E=0.016; F=0.013; G=0.016; H=0.015 },
    @{ Row=39; N="500";   C=0.013; D=0.016; E=0.016; F=0.01;  G=0.046; H=0.004 },
    @{ Row=40; N="1000";  C=0.003; D=0.007; E=0.007; F=0.003; G=0.008999999999999999; H=0.013 },
    @{ Row=41; N="2500";  C=0.004; D=0.005; E=0.005; F=0.005; G=0.011; H=0.008999999999999999 },
    @{ Row=42; N="5000";  C=0.005; D=0.006; E=0.002; F=0.002; G=0.003; H=0.005 },
    @{ Row=43; N="10000"; C=0.004; D=0.003; E=0.001; F=0.001; G=0.002; H=0.003 }
)

foreach ($r in $accBalStdRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "accuracy_balanced_std"
    $ws.Range("B$row").Value = "'" + $r.N
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
}

# ---------------------------------------------------------------------
# 4) Match the bold/border/centered style used for the other metric-name
#    cells in column A (style carried by e.g. A2) on our new A-column
#    cells, and clear the leftover formatting the row-insert copied down.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A16:A22").PasteSpecial(-4122)
$ws.Range("A37:A43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "Added accuracy_balanced_mean (rows 16-22) and accuracy_balanced_std (rows 37-43)."
